# Apply targeted odds updates to Sheet1, matching the committed diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5

# Row 3 updates
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.6
$ws.Range("J3").Value = 3.5
$ws.Range("W3").Value = 8
$ws.Range("Y3").Value = 11
$ws.Range("AD3").Value = 6
$ws.Range("AN3").Value = 4.75
$ws.Range("AO3").Value = 17
$ws.Range("AP3").Value = 29
